# Auto-generated edit script applying profit-sheet value updates
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 40.5
$ws.Range("I11").Value = 40.5
$ws.Range("K11").Value = 40.5
$ws.Range("M11").Value = 99.5
$ws.Range("H51").Value = 34980.832
$ws.Range("I51").Value = 42498.5
$ws.Range("K51").Value = 42498.5
$ws.Range("M51").Value = -42014.5
$ws.Range("H75").Value = 43000
$ws.Range("J75").Value = 43000
$ws.Range("L75").Value = 43000
$ws.Range("N75").Value = -44872
$ws.Range("H78").Value = 43000
$ws.Range("J78").Value = 43000
$ws.Range("L78").Value = 129000
$ws.Range("N78").Value = -138360
$ws.Range("H98").Value = 47622110
$ws.Range("I98").Value = 50003170
$ws.Range("K98").Value = 50003170
$ws.Range("M98").Value = -50001672
$ws.Range("H105").Value = 54381
$ws.Range("J105").Value = 54381
$ws.Range("L105").Value = 54381
$ws.Range("N105").Value = -61369
$ws.Range("H122").Value = 47622110
$ws.Range("I122").Value = 50003170
$ws.Range("K122").Value = 150009510
$ws.Range("M122").Value = -150007060
$ws.Range("H137").Value = 4708.3687
$ws.Range("I137").Value = 3373.4167
$ws.Range("J137").Value = 6996.857
$ws.Range("K137").Value = 10120.2501
$ws.Range("L137").Value = 20990.571
$ws.Range("M137").Value = -7570.250100000001
$ws.Range("N137").Value = -26090.571

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2258.7273
$ws.Range("J45").Value = 3014
$ws.Range("L45").Value = 3014
$ws.Range("N45").Value = -3768
$ws.Range("H74").Value = 22641.49
$ws.Range("I74").Value = 26335.176
$ws.Range("K74").Value = 26335.176
$ws.Range("M74").Value = -25461.176
$ws.Range("H77").Value = 22641.49
$ws.Range("I77").Value = 26335.176
$ws.Range("K77").Value = 131675.88
$ws.Range("M77").Value = -127307.88
$ws.Range("H122").Value = 3687.9412
$ws.Range("I122").Value = 2093.0476
$ws.Range("K122").Value = 6279.1428
$ws.Range("M122").Value = -3829.1428

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2155.125
$ws.Range("I105").Value = 2082.0476
$ws.Range("J105").Value = 2666.6667
$ws.Range("K105").Value = 2082.0476
$ws.Range("L105").Value = 2666.6667
$ws.Range("M105").Value = -335.0475999999999
$ws.Range("N105").Value = -6160.6667
$ws.Range("H129").Value = 500025000
$ws.Range("J129").Value = 1000000000
$ws.Range("L129").Value = 1000000000
$ws.Range("N129").Value = -1000010000

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3681.2432
$ws.Range("I16").Value = 1368.8572
$ws.Range("J16").Value = 5088.7827
$ws.Range("K16").Value = 1368.8572
$ws.Range("L16").Value = 5088.7827
$ws.Range("M16").Value = -1081.8572
$ws.Range("N16").Value = -5662.7827
$ws.Range("H31").Value = 4830.2354
$ws.Range("I31").Value = 2433.647
$ws.Range("K31").Value = 2433.647
$ws.Range("M31").Value = -2138.647
$ws.Range("H34").Value = 4830.2354
$ws.Range("I34").Value = 2433.647
$ws.Range("K34").Value = 2433.647
$ws.Range("M34").Value = -2231.647
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
$ws.Range("H86").Value = 28423374
$ws.Range("I86").Value = 13490980
$ws.Range("K86").Value = 13490980
$ws.Range("M86").Value = -13489857
$ws.Range("H89").Value = 28423374
$ws.Range("I89").Value = 13490980
$ws.Range("K89").Value = 67454900
$ws.Range("M89").Value = -67449284
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()
$ws.Range("H113").Value = 3681.2432
$ws.Range("I113").Value = 1368.8572
$ws.Range("J113").Value = 5088.7827
$ws.Range("K113").Value = 1368.8572
$ws.Range("L113").Value = 5088.7827
$ws.Range("M113").Value = 801.1428000000001
$ws.Range("N113").Value = -9428.7827
$ws.Range("H134").Value = 3984.3547
$ws.Range("I134").Value = 2197.2632
$ws.Range("K134").Value = 6591.7896
$ws.Range("M134").Value = -4056.7896
$ws.Range("H141").Value = 373332.5
$ws.Range("J141").Value = 373332.5
$ws.Range("L141").Value = 373332.5
$ws.Range("N141").Value = -383692.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 650
$ws.Range("I3").Value = 650
$ws.Range("K3").Value = 1950
$ws.Range("M3").Value = -1838
$ws.Range("H11").Value = 1407.5555
$ws.Range("I11").Value = 896
$ws.Range("K11").Value = 2688
$ws.Range("M11").Value = -2548
$ws.Range("H134").Value = 45335.36
$ws.Range("I134").Value = 59666.723
$ws.Range("K134").Value = 179000.169
$ws.Range("M134").Value = -173930.169

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H104").Value = 87899.5
$ws.Range("J104").Value = 87899.5
$ws.Range("L104").Value = 87899.5
$ws.Range("N104").Value = -94887.5
$ws.Range("H107").Value = 571810.5
$ws.Range("J107").Value = 274.83334
$ws.Range("L107").Value = 274.83334
$ws.Range("N107").Value = -4114.83334
$ws.Range("H122").Value = 2590003.8
$ws.Range("I122").Value = 3152060.5
$ws.Range("J122").Value = 4542.4
$ws.Range("K122").Value = 9456181.5
$ws.Range("L122").Value = 13627.2
$ws.Range("M122").Value = -9453731.5
$ws.Range("N122").Value = -18527.2
$ws.Range("H126").Value = 19234104
$ws.Range("J126").Value = 5069.6665
$ws.Range("L126").Value = 15208.9995
$ws.Range("N126").Value = -20148.9995
$ws.Range("H139").Value = 66664
$ws.Range("J139").Value = 66664
$ws.Range("L139").Value = 66664
$ws.Range("N139").Value = -76944

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5111.8335
$ws.Range("I7").Value = 4252.385
$ws.Range("J7").Value = 6127.5454
$ws.Range("K7").Value = 4252.385
$ws.Range("L7").Value = 6127.5454
$ws.Range("M7").Value = -4140.385
$ws.Range("N7").Value = -6351.5454
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("H126").Value = 5111.8335
$ws.Range("I126").Value = 4252.385
$ws.Range("J126").Value = 6127.5454
$ws.Range("K126").Value = 12757.155
$ws.Range("L126").Value = 18382.6362
$ws.Range("M126").Value = -10287.155
$ws.Range("N126").Value = -23322.6362
$ws.Range("H132").Value = 7252678
$ws.Range("I132").Value = 12197455
$ws.Range("K132").Value = 36592365
$ws.Range("M132").Value = -36589835

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H136").Value = 19426134
$ws.Range("J136").Value = 460407.47
$ws.Range("L136").Value = 1381222.41
$ws.Range("N136").Value = -1386322.41
